$wb = $excel.ActiveWorkbook

# Rename sheets: HU_High -> PL_High, HU_Low -> PL_Low
$wb.Worksheets.Item("HU_High").Name = "PL_High"
$wb.Worksheets.Item("HU_Low").Name = "PL_Low"
